$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# H21: clarify telecentro/PRFV observation
$ws.Range("H21").Value = "Telecentro ya traspaso su nodo solo falta desmontar el prfv que quedo"

# H66: clarify "Picada" observation
$ws.Range("H66").Value = "fue bajada como que no es de la empresa pero no se ve que sea de telecentro reverificar o esperar a Pedro para verla"

# Row 80: case number and OT updated (Caso is stored as text, so force text
# to avoid Excel auto-converting the numeric-looking string to a number,
# then restore the default "Normal" style so no stray formatting is left)
$ws.Range("A80").NumberFormat = "@"
$ws.Range("A80").Value = "6578"
$ws.Range("A80").Style = "Normal"
$ws.Range("E80").Value = "Pendiente ADM"

# Old row 83 (Caso -545, Jeronimo Salguero 3601) is removed entirely; rows below shift up
$ws.Rows("83").Delete()
